$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns in this sheet store plain text
# (percent-like and dotted-thousands price strings), not real numbers/percentages.
# Guard the whole data range with a Text number format while writing the new
# values so Excel's COM layer doesn't silently reinterpret strings like
# "1.010" / "0.08500" / "0.00001073" as numbers (which would drop trailing
# zeros or switch to scientific notation). Afterwards restore the cells to
# the workbook's default "Normal" style so no stray formatting is introduced.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.486.70"
$ws.Range("D3").Value = "1.957.61"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "322.29"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("D8").Value = "0.4071"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").Value = "53.90"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").Value = "0.08500"
$ws.Range("E10").Value = "  -7.38%  "
$ws.Range("D11").Value = "1.059"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "1.946.79"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "7.586"
$ws.Range("E14").Value = "  -3.48%  "
$ws.Range("D15").Value = "6.168"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "90.51"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "0.00001073"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "0.06630"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "18.51"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "5.849"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "28.483.76"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").Value = "11.43"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "2.160.69"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "20.32"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "2.174"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").Value = "5.839"
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("D31").Value = "124.36"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "0.9855"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").Value = "0.09658"
$ws.Range("D34").Value = "1.454"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").Value = "5.640"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").Value = "3.690"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "9.107"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "0.02335"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("D39").Value = "0.06194"
$ws.Range("D40").Value = "1.253"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").Value = "0.6238"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").Value = "11.20"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "0.1918"
$ws.Range("E44").Value = "  -3.37%  "
$ws.Range("D45").Value = "1.355"
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("D46").Value = "0.5956"
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("D47").Value = "13.02"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "2.062"
$ws.Range("E48").Value = "  -5.25%  "
$ws.Range("D49").Value = "3.410"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "0.06811"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "111.19"
$ws.Range("E51").Value = "  -1.17%  "

$dataRange.Style = "Normal"
